# comment from python script
$wb = $excel.ActiveWorkbook

# --- Sheet: Astronauta (sheet1) ---
$ws = $wb.Worksheets.Item("Astronauta")
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("B20").Select()

# --- Sheet: Senador (sheet2) ---
$ws = $wb.Worksheets.Item("Senador")
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("B5").Select()

# --- Sheet: Mago (sheet3) ---
$ws = $wb.Worksheets.Item("Mago")
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 1
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("B22").Select()

# --- Sheet: Ninja (sheet4) ---
$ws = $wb.Worksheets.Item("Ninja")
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("C4").Select()

$ws.Activate()
